$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "two_line": append rows 109-113
# Columns: A detected_date(date) B breakout_date C Time_Frame D stockname
#          E date1(date) F value1 G date2(date) H value2 I buyORsell J "Date Time"
# ---------------------------------------------------------------------------
$wsTwoLine = $wb.Worksheets.Item("two_line")

$twoLineRows = @(
    @{ A = 45443.42708333334; B = "14-06-2024 10:15:00"; C = "hour"; D = "ASIANPAINT.NS"; E = 45435.38541666666; F = 2921.25;             G = 45441.55208333334; H = 2921;                 I = "High"; J = "14/06/2024 04:48:06" },
    @{ A = 45455.59375;       B = "14-06-2024 10:15:00"; C = "hour"; D = "VOLTAS.NS";     E = 45449.38541666666; F = 1477.849975585938;    G = 45454.42708333334; H = 1476;                 I = "High"; J = "14/06/2024 04:48:06" },
    @{ A = 45427.42708333334; B = "14-06-2024 09:15:00"; C = "hour"; D = "KPEL.BO";       E = 45425.46875;       F = 422.1000061035156;    G = 45425.55208333334; H = 422.1000061035156;     I = "Low";  J = "14/06/2024 04:48:06" },
    @{ A = 45434.51041666666; B = "14-06-2024 10:15:00"; C = "hour"; D = "NIITLTD.NS";    E = 45429.51041666666; F = 104.9499969482422;    G = 45429.63541666666; H = 104.9499969482422;     I = "High"; J = "14/06/2024 04:48:06" },
    @{ A = 45447.55208333334; B = "14-06-2024 09:15:00"; C = "hour"; D = "ADFFOODS.NS";   E = 45441.42708333334; F = 226.8000030517578;    G = 45446.38541666666; H = 227.25;                I = "High"; J = "14/06/2024 04:48:06" }
)

$r = 109
foreach ($row in $twoLineRows) {
    $wsTwoLine.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsTwoLine.Cells.Item($r, 1).Value2 = $row.A
    $wsTwoLine.Cells.Item($r, 2).Value2 = $row.B
    $wsTwoLine.Cells.Item($r, 3).Value2 = $row.C
    $wsTwoLine.Cells.Item($r, 4).Value2 = $row.D
    $wsTwoLine.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsTwoLine.Cells.Item($r, 5).Value2 = $row.E
    $wsTwoLine.Cells.Item($r, 6).Value2 = $row.F
    $wsTwoLine.Cells.Item($r, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsTwoLine.Cells.Item($r, 7).Value2 = $row.G
    $wsTwoLine.Cells.Item($r, 8).Value2 = $row.H
    $wsTwoLine.Cells.Item($r, 9).Value2 = $row.I
    $wsTwoLine.Cells.Item($r, 10).Value2 = $row.J
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "ph_pl_breakout_line": append rows 938-947
# Columns: A stockname B Datetime(date) C High D Low E Close F isPivot
#          G PHorPLValue H time_frame I TdyDate J TdyClose K PClose L "Date Time"
# ---------------------------------------------------------------------------
$wsBreakout = $wb.Worksheets.Item("ph_pl_breakout_line")

$breakoutRows = @(
    @{ A = "RIR.BO";       B = 45442.38541666666; C = 1789.800048828125;    D = 1700;                 E = 1710.050048828125;    F = "High"; G = 1789.800048828125;    H = "hour"; I = "14-06-2024 09:15:00"; J = 1823.25;             K = 1736.449951171875; L = "14/06/2024 04:48:06" },
    @{ A = "RIR.BO";       B = 45446.38541666666; C = 1782.800048828125;    D = 1655.5;               E = 1678;                 F = "High"; G = 1782.800048828125;    H = "hour"; I = "14-06-2024 09:15:00"; J = 1823.25;             K = 1736.449951171875; L = "14/06/2024 04:48:06" },
    @{ A = "RIR.BO";       B = 45454.38541666666; C = 1760;                 D = 1690;                 E = 1699;                 F = "High"; G = 1760;                 H = "hour"; I = "14-06-2024 09:15:00"; J = 1823.25;             K = 1736.449951171875; L = "14/06/2024 04:48:06" },
    @{ A = "MOSCHIP.BO";   B = 45453.38541666666; C = 197.3000030517578;    D = 188;                  E = 189.6000061035156;    F = "High"; G = 197.3000030517578;    H = "hour"; I = "14-06-2024 09:15:00"; J = 209.9499969482422;    K = 197;               L = "14/06/2024 04:48:06" },
    @{ A = "GPTINFRA.NS";  B = 45450.38541666666; C = 259.7999877929688;    D = 247.3000030517578;    E = 255.25;               F = "High"; G = 259.7999877929688;    H = "hour"; I = "14-06-2024 10:15:00"; J = 260;                 K = 259.7999877929688; L = "14/06/2024 04:48:06" },
    @{ A = "GOODYEAR.BO";  B = 45454.38541666666; C = 1179.75;              D = 1156.050048828125;    E = 1170.050048828125;    F = "High"; G = 1179.75;              H = "hour"; I = "14-06-2024 09:15:00"; J = 1182.599975585938;    K = 1179.599975585938; L = "14/06/2024 04:48:06" },
    @{ A = "SANJIVIN.BO";  B = 45454.38541666666; C = 176.8000030517578;    D = 170.1999969482422;    E = 170.5500030517578;    F = "High"; G = 176.8000030517578;    H = "hour"; I = "14-06-2024 09:15:00"; J = 177.8500061035156;    K = 175.8999938964844; L = "14/06/2024 04:48:06" },
    @{ A = "SKYGOLD.NS";   B = 45453.38541666666; C = 1297.949951171875;    D = 1269;                 E = 1269.949951171875;    F = "High"; G = 1297.949951171875;    H = "hour"; I = "14-06-2024 09:15:00"; J = 1337.400024414062;    K = 1273.75;           L = "14/06/2024 04:48:06" },
    @{ A = "UNIAUTO.BO";   B = 45443.38541666666; C = 171.9499969482422;    D = 164;                  E = 167.4499969482422;    F = "High"; G = 171.9499969482422;    H = "hour"; I = "14-06-2024 09:15:00"; J = 172.9499969482422;    K = 170.1999969482422; L = "14/06/2024 04:48:06" },
    @{ A = "UNIAUTO.BO";   B = 45454.38541666666; C = 171.9499969482422;    D = 163.1999969482422;    E = 166.5;                F = "High"; G = 171.9499969482422;    H = "hour"; I = "14-06-2024 09:15:00"; J = 172.9499969482422;    K = 170.1999969482422; L = "14/06/2024 04:48:06" }
)

$r = 938
foreach ($row in $breakoutRows) {
    $wsBreakout.Cells.Item($r, 1).Value2 = $row.A
    $wsBreakout.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsBreakout.Cells.Item($r, 2).Value2 = $row.B
    $wsBreakout.Cells.Item($r, 3).Value2 = $row.C
    $wsBreakout.Cells.Item($r, 4).Value2 = $row.D
    $wsBreakout.Cells.Item($r, 5).Value2 = $row.E
    $wsBreakout.Cells.Item($r, 6).Value2 = $row.F
    $wsBreakout.Cells.Item($r, 7).Value2 = $row.G
    $wsBreakout.Cells.Item($r, 8).Value2 = $row.H
    $wsBreakout.Cells.Item($r, 9).Value2 = $row.I
    $wsBreakout.Cells.Item($r, 10).Value2 = $row.J
    $wsBreakout.Cells.Item($r, 11).Value2 = $row.K
    $wsBreakout.Cells.Item($r, 12).Value2 = $row.L
    $r++
}
